# Applies the "Added workflow image_netayume_lumina_t2i_mod.json Anime" edit
# to the first ("Workflows") table on the page.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Resize the three columns (3295/3028/3027 dxa -> 4135/2585/2630 dxa). ---
# Word's Column.Width is expressed in points, i.e. dxa / 20.
$t.Columns.Item(1).Width = 4135 / 20
$t.Columns.Item(2).Width = 2585 / 20
$t.Columns.Item(3).Width = 2630 / 20

# --- Fill in the previously-empty Status / Notes cells for existing rows. ---

# Row 4: flux_quick.json
$t.Cell(4, 2).Range.Text = "Work, some image errors"
$t.Cell(4, 3).Range.Text = "Images are mostly good, ~5sec"

# Row 5: SDXLturbo_Quick2.json
$t.Cell(5, 2).Range.Text = "Work, some image errors, distortions in hands and faces. "
$t.Cell(5, 3).Range.Text = "Images mostly good, Fast, 0.5sec per image"

# Row 6: Qwen Image Rapid.json
$t.Cell(6, 2).Range.Text = "Works, rare image errors"
$t.Cell(6, 3).Range.Text = "Image look good, ~9sec"

# --- Insert a brand-new row for the newly added workflow, right before the
#     trailing blank row (so it becomes row 7, and the blank row stays last).
$blankRow = $t.Rows.Item(7)
$newRow = $t.Rows.Add($blankRow)
$newRow.Cells.Item(1).Range.Text = "image_netayume_lumina_t2i_ex.json"

Write-Output "Workflows table updated: $($t.Rows.Count) rows."
